# Auto-generated script to apply Goblin_Profits scheduled-runner update
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ALC_updates = @{
    "H64" = 7335.2646
    "I64" = 3666.6667
    "J64" = 8121.393
    "K64" = 3666.6667
    "L64" = 8121.393
    "M64" = -3418.6667
    "N64" = -8617.393
    "H67" = 7335.2646
    "I67" = 3666.6667
    "J67" = 8121.393
    "K67" = 3666.6667
    "L67" = 8121.393
    "M67" = -2808.6667
    "N67" = -9837.393
    "H74" = 14424.272
    "I74" = 14424.272
    "K74" = 14424.272
    "M74" = -13488.272
    "H77" = 14424.272
    "I77" = 14424.272
    "K77" = 72121.36
    "M77" = -67441.36
    "H100" = 5587.706
    "I100" = 2284.7144
    "K100" = 2284.7144
    "M100" = -1743.7144
    "H137" = 1498.5555
    "I137" = 1518
    "K137" = 4554
    "M137" = -2004
}
foreach ($cellRef in $ALC_updates.Keys) {
    $ws.Range($cellRef).Value = $ALC_updates[$cellRef]
}

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ARM_updates = @{
    "H2" = 2471.8823
    "I2" = 989
    "K2" = 989
    "M2" = -876
    "H74" = 1688.2307
    "I74" = 1783.6666
    "K74" = 1783.6666
    "M74" = -909.6666
    "H77" = 1688.2307
    "I77" = 1783.6666
    "K77" = 8918.333000000001
    "M77" = -4550.333000000001
    "H116" = 2471.8823
    "I116" = 989
    "K116" = 989
    "M116" = 1305
    "H132" = 2095.9788
    "I132" = 2084.6047
    "J132" = 2218.25
    "K132" = 6253.8141
    "L132" = 6654.75
    "M132" = -3723.8141
    "N132" = -11714.75
}
foreach ($cellRef in $ARM_updates.Keys) {
    $ws.Range($cellRef).Value = $ARM_updates[$cellRef]
}

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$BSM_updates = @{
    "H3" = 2471.8823
    "I3" = 989
    "K3" = 989
    "M3" = -875
    "H82" = 26037.4
    "J82" = 60000
    "L82" = 60000
    "N82" = -60766
    "H85" = 26037.4
    "J85" = 60000
    "L85" = 60000
    "N85" = -62652
    "H94" = 2992.6667
    "I94" = 1852.6111
    "K94" = 1852.6111
    "M94" = -1401.6111
}
foreach ($cellRef in $BSM_updates.Keys) {
    $ws.Range($cellRef).Value = $BSM_updates[$cellRef]
}

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$CRP_updates = @{
    "H31" = 4402.1577
    "I31" = 1805.2222
    "K31" = 1805.2222
    "M31" = -1510.2222
    "H34" = 4402.1577
    "I34" = 1805.2222
    "K34" = 1805.2222
    "M34" = -1603.2222
    "H58" = 1191.1818
    "I58" = 1121
    "K58" = 1121
    "M58" = -918
    "H136" = 1191.1818
    "I136" = 1121
    "K136" = 3363
    "M136" = -813
}
foreach ($cellRef in $CRP_updates.Keys) {
    $ws.Range($cellRef).Value = $CRP_updates[$cellRef]
}

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$CUL_updates = @{
    "H87" = 5318.6665
    "I87" = 5318.6665
    "J87" = 0
    "K87" = 15955.9995
    "L87" = 0
    "M87" = -14707.9995
    "H90" = 5318.6665
    "I90" = 5318.6665
    "J90" = 0
    "K90" = 47867.9985
    "L90" = 0
    "M90" = -41627.9985
    "H134" = 4261.9443
    "I134" = 1530.7142
    "K134" = 4592.142599999999
    "M134" = 477.8574000000008
}
foreach ($cellRef in $CUL_updates.Keys) {
    $ws.Range($cellRef).Value = $CUL_updates[$cellRef]
}
$ws.Range("N87").ClearContents()
$ws.Range("N90").ClearContents()

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$GSM_updates = @{
    "H97" = 754.4
    "I97" = 734.75
    "K97" = 734.75
    "M97" = -238.75
    "H132" = 1988.9744
    "I132" = 1988.9744
    "K132" = 5966.9232
    "M132" = -3436.9232
}
foreach ($cellRef in $GSM_updates.Keys) {
    $ws.Range($cellRef).Value = $GSM_updates[$cellRef]
}

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$LTW_updates = @{
    "H7" = 6615.8
    "I7" = 4711.625
    "K7" = 4711.625
    "M7" = -4599.625
    "H22" = 3548.9644
    "I22" = 3610.6365
    "J22" = 3509.0588
    "K22" = 3610.6365
    "L22" = 3509.0588
    "M22" = -3315.6365
    "N22" = -4099.0588
    "H27" = 3548.9644
    "I27" = 3610.6365
    "J27" = 3509.0588
    "K27" = 3610.6365
    "L27" = 3509.0588
    "M27" = -3503.6365
    "N27" = -3723.0588
    "H74" = 44125
    "I74" = 38833.332
    "K74" = 38833.332
    "M74" = -37835.332
    "H77" = 44125
    "I77" = 38833.332
    "K77" = 116499.996
    "M77" = -111507.996
    "H126" = 6615.8
    "I126" = 4711.625
    "K126" = 14134.875
    "M126" = -11664.875
    "H132" = 3510.4412
    "I132" = 3296.423
    "K132" = 9889.269
    "M132" = -7359.269
    "H136" = 12236.827
    "I136" = 1759.8889
    "J136" = 14429.675
    "K136" = 5279.6667
    "L136" = 43289.02499999999
    "M136" = -2729.6667
    "N136" = -48389.02499999999
}
foreach ($cellRef in $LTW_updates.Keys) {
    $ws.Range($cellRef).Value = $LTW_updates[$cellRef]
}

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$WVR_updates = @{
    "H132" = 1428.5807
    "I132" = 1435.6
    "J132" = 1415.8182
    "K132" = 4306.799999999999
    "L132" = 4247.4546
    "M132" = -1776.799999999999
    "N132" = -9307.454600000001
    "H136" = 920.1
    "I136" = 921.1786
    "J136" = 905
    "K136" = 2763.5358
    "L136" = 2715
    "M136" = -213.5357999999997
    "N136" = -7815
}
foreach ($cellRef in $WVR_updates.Keys) {
    $ws.Range($cellRef).Value = $WVR_updates[$cellRef]
}
